$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of (x, y, z) values inserted at row 2, followed by the original
# data rows 2-20 shifted down by one position each. (Original rows 21 and
# 22 are dropped, matching the new dimension A1:C21.)
$data = @(
  @(-0.0114537235349416, -0.0215329993516206, 0.0035124751739203),
  @(-0.0171042270958423, -0.0526871271431446, 0.0519235469400882),
  @(0.0067195175215601, 0.1915062516927719, -0.0238237436860799),
  @(-0.001527163083665, 0.3747658133506775, 0.0192422550171613),
  @(-0.3475823104381561, 0.5015203952789307, -0.2428189367055893),
  @(-0.2506074607372284, 0.8736900091171265, -0.1963931769132614),
  @(-0.168751522898674, 0.936914563179016, -0.131183311343193),
  @(-0.2408336251974105, -0.064446285367012, -0.3248275816440582),
  @(-0.2364048510789871, 0.18539759516716, -1.519680023193359),
  @(-1.216843605041504, -1.400255799293518, -1.343445420265198),
  @(1.094823241233826, -2.151772737503052, 0.7765624523162842),
  @(-0.2055561542510986, -1.607033729553223, -0.9816604256629944),
  @(-0.5111414790153503, -1.858557462692261, 0.0117591563612222),
  @(-0.3014619946479797, -0.9010262489318848, 0.4983133375644684),
  @(-0.2797762751579284, -1.265254616737366, 0.2492330223321914),
  @(0.2038762718439102, -1.481042742729187, -0.2376265823841095),
  @(-0.328340083360672, 0.6785185933113098, 0.7831292152404785),
  @(0.0797179117798805, 0.1185078546404838, 0.2952006161212921),
  @(0.0479529201984405, 0.3949243724346161, -0.1617265790700912),
  @(-0.08704829961061469, -0.058643065392971, -0.2785545587539673)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $data[$i][0]
  $ws.Cells.Item($row, 2).Value = $data[$i][1]
  $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Drop the now-surplus last row (previously row 22) so the sheet ends at row 21
$ws.Rows.Item(22).Delete()
